# Update the "想去人数" (interested-count) values in column F on both the
# "展览" and "全部类型" worksheets to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    "F2"  = 15210
    "F3"  = 19634
    "F5"  = 186
    "F14" = 229
    "F15" = 260
    "F17" = 1540
    "F18" = 150
    "F20" = 122
    "F21" = 252
    "F22" = 8295
    "F23" = 995
    "F26" = 73
    "F27" = 1283
    "F28" = 37
    "F29" = 15
    "F31" = 6634
    "F32" = 141
    "F34" = 194
    "F37" = 5679
}
foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    "F2"  = 15210
    "F3"  = 19634
    "F5"  = 186
    "F14" = 229
    "F15" = 260
    "F17" = 1540
    "F18" = 150
    "F21" = 122
    "F22" = 252
    "F23" = 8295
    "F24" = 995
    "F27" = 73
    "F28" = 1283
    "F29" = 37
    "F30" = 15
    "F34" = 6634
    "F35" = 141
    "F37" = 194
    "F40" = 5679
}
foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
